# "update to latex document"
#
# The sheet gained three new helper columns (N, O, P) next to the existing
# data table, column L got a bit wider to fit them, the result chart was
# repositioned/resized to make room, and the window/selection was left
# where the author was last working (C32, zoomed in to 130%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new percentage-formatted helper cells (N2:P2, N4:P4) -----------------
# N2 only carries the percentage number format (no value); O2/P2/N4/O4/P4
# carry plain numeric fractions that Excel renders with the General format.
$ws.Range("N2").NumberFormat = "0%"
$ws.Range("N2").Value = $null

$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 0.7

$ws.Range("N4").Value = 0.4
$ws.Range("O4").Value = 0.6
$ws.Range("P4").Value = 0.8

# --- widen the new column L (index 12) -------------------------------------
$ws.Columns(12).ColumnWidth = 18.14

# --- move/resize the results chart to sit over the (now wider) table ------
$co = $ws.ChartObjects(1)
$co.Left = 275.23828125
$co.Top = 262.87496062992125
$co.Width = 522.6611328125
$co.Height = 229.87503937007875

# --- window/view state: zoomed to 130%, selection parked on C32 -----------
$win = $excel.ActiveWindow
$win.Zoom = 130
$ws.Range("C32").Select() | Out-Null
